$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the now-extra data rows (old rows 8, 9, 10) since the new
# TPM data only has 6 data rows (rows 2-7) instead of 9.
$ws.Rows.Item(8).Delete()
$ws.Rows.Item(8).Delete()
$ws.Rows.Item(8).Delete()

# Update remaining data rows (2-7) with the refreshed TPM-derived values.
$ws.Range("A2").Value = "FAPs"
$ws.Range("B2").Value = "Agt"
$ws.Range("C2").Value = "Agtr1a"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 2
$ws.Range("F2").Value = 0.6666666666666666
$ws.Range("G2").Value = 1.551523333333333
$ws.Range("H2").Value = 4.65457
$ws.Range("I2").Value = 0.950716861801202
$ws.Range("J2").Value = 0.950716861801202
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 1.137783333333333
$ws.Range("N2").Value = 3.41335
$ws.Range("O2").Value = 0.02543039699931523
$ws.Range("P2").Value = 0.02543039699931523
$ws.Range("Q2").Value = 1.765297389944444
$ws.Range("R2").Value = 15.8876765095
$ws.Range("S2").Value = 0.02417710722954768
$ws.Range("T2").Value = 0.02417710722954768

$ws.Range("A3").Value = "FAPs"
$ws.Range("B3").Value = "Agt"
$ws.Range("C3").Value = "Agtr1a"
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 2
$ws.Range("F3").Value = 0.6666666666666666
$ws.Range("G3").Value = 1.551523333333333
$ws.Range("H3").Value = 4.65457
$ws.Range("I3").Value = 0.950716861801202
$ws.Range("J3").Value = 0.950716861801202
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 15.44817966666667
$ws.Range("N3").Value = 46.344539
$ws.Range("O3").Value = 0.3452795715412271
$ws.Range("P3").Value = 0.345279571541227
$ws.Range("Q3").Value = 23.96821121035888
$ws.Range("R3").Value = 215.71390089323
$ws.Range("S3").Value = 0.328263110699739
$ws.Range("T3").Value = 0.328263110699739

$ws.Range("A4").Value = "FAPs"
$ws.Range("B4").Value = "Agt"
$ws.Range("C4").Value = "Agtr1a"
$ws.Range("D4").Value = "MuSCs"
$ws.Range("E4").Value = 2
$ws.Range("F4").Value = 0.6666666666666666
$ws.Range("G4").Value = 1.551523333333333
$ws.Range("H4").Value = 4.65457
$ws.Range("I4").Value = 0.950716861801202
$ws.Range("J4").Value = 0.950716861801202
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 28.15511333333333
$ws.Range("N4").Value = 84.46534
$ws.Range("O4").Value = 0.6292900314594577
$ws.Range("P4").Value = 0.6292900314594577
$ws.Range("Q4").Value = 43.6833152893111
$ws.Range("R4").Value = 393.1498376037999
$ws.Range("S4").Value = 0.5982766438719153
$ws.Range("T4").Value = 0.5982766438719153

$ws.Range("A5").Value = "MuSCs"
$ws.Range("B5").Value = "Agt"
$ws.Range("C5").Value = "Agtr1a"
$ws.Range("D5").Value = "ECs"
$ws.Range("E5").Value = 1
$ws.Range("F5").Value = 0.3333333333333333
$ws.Range("G5").Value = 0.08042766666666666
$ws.Range("H5").Value = 0.241283
$ws.Range("I5").Value = 0.04928313819879805
$ws.Range("J5").Value = 0.04928313819879804
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 1.137783333333333
$ws.Range("N5").Value = 3.41335
$ws.Range("O5").Value = 0.02543039699931523
$ws.Range("P5").Value = 0.02543039699931523
$ws.Range("Q5").Value = 0.09150925867222222
$ws.Range("R5").Value = 0.8235833280500001
$ws.Range("S5").Value = 0.001253289769767551
$ws.Range("T5").Value = 0.001253289769767551

$ws.Range("A6").Value = "MuSCs"
$ws.Range("B6").Value = "Agt"
$ws.Range("C6").Value = "Agtr1a"
$ws.Range("D6").Value = "FAPs"
$ws.Range("E6").Value = 1
$ws.Range("F6").Value = 0.3333333333333333
$ws.Range("G6").Value = 0.08042766666666666
$ws.Range("H6").Value = 0.241283
$ws.Range("I6").Value = 0.04928313819879805
$ws.Range("J6").Value = 0.04928313819879804
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 15.44817966666667
$ws.Range("N6").Value = 46.344539
$ws.Range("O6").Value = 0.3452795715412271
$ws.Range("P6").Value = 0.345279571541227
$ws.Range("Q6").Value = 1.242461044837444
$ws.Range("R6").Value = 11.182149403537
$ws.Range("S6").Value = 0.01701646084148807
$ws.Range("T6").Value = 0.01701646084148807

$ws.Range("A7").Value = "MuSCs"
$ws.Range("B7").Value = "Agt"
$ws.Range("C7").Value = "Agtr1a"
$ws.Range("D7").Value = "MuSCs"
$ws.Range("E7").Value = 1
$ws.Range("F7").Value = 0.3333333333333333
$ws.Range("G7").Value = 0.08042766666666666
$ws.Range("H7").Value = 0.241283
$ws.Range("I7").Value = 0.04928313819879805
$ws.Range("J7").Value = 0.04928313819879804
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 28.15511333333333
$ws.Range("N7").Value = 84.46534
$ws.Range("O7").Value = 0.6292900314594577
$ws.Range("P7").Value = 0.6292900314594577
$ws.Range("Q7").Value = 2.264450070135555
$ws.Range("R7").Value = 20.38005063122
$ws.Range("S7").Value = 0.03101338758754243
$ws.Range("T7").Value = 0.03101338758754242

